# Admin pages, get method
# Adds the "video demonstrating source code" hyperlink (B10) next to the
# existing GitHub repository hyperlink (B11), mirroring its formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$videoUrl = "https://youtu.be/JILqNlTV-8k"

# 1. Put the display text in B10 (shared string, same text as the URL).
$ws.Range("B10").Value = $videoUrl

# 2. Turn B10 into a real hyperlink pointing at the YouTube video.
$ws.Hyperlinks.Add($ws.Range("B10"), $videoUrl) | Out-Null

# 3. Copy B11's formatting (font/fill/alignment - the "Hiperlink" look) onto
#    B10 so both rows look the same, then clear the clipboard marquee.
$ws.Range("B11").Copy() | Out-Null
$ws.Range("B10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# 4. Match the saved selection (B10 was the last selected cell).
$ws.Range("B10").Select()
